$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "d2" column (old column B) was a duplicate of column A ("d1") and is
# being dropped from the exported CSV/metadata. Deleting the whole column
# shifts the old C ("l1") and D ("l2") columns left and keeps all of their
# values and number formatting intact.
$ws.Columns("B").Delete()

# Rename the remaining first column header from "d1" to "depth".
$ws.Range("A1").Value = "depth"

# Re-apply the sort so the sort state / sort conditions only reference the
# three remaining columns (A:C) instead of the old four-column range.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A31"))
$ws.Sort.SortFields.Add($ws.Range("B2:B31"))
$ws.Sort.SortFields.Add($ws.Range("C2:C31"))
$ws.Sort.SetRange($ws.Range("A1:C31"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Leave the selection where the editor ended up after the edit.
[void]$ws.Range("D2").Select()
